# Generate Report for Archive
#
# The report rows for the three outstanding files (1b48e175, 3f14bde9,
# ba0dd1f4) are re-sorted: 3f14bde9 and ba0dd1f4 move up to rows 3 and 4
# (their status flips from "Ready for handoff" to "In Translation" and
# their handoff timestamps refresh), while 1b48e175 drops to row 5
# (still "Ready for handoff").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Row 3 -> 3f14bde9
$ws1.Range("A3").Value = "3f14bde9-0c56-4a70-af96-9a02074d8301.md"
$ws1.Range("B3").Value = "e2e\3f14bde9-0c56-4a70-af96-9a02074d8301.md"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"
$ws1.Range("G3").Value = "2016-08-29 22:42:50"

# Row 4 -> ba0dd1f4
$ws1.Range("A4").Value = "ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md"
$ws1.Range("B4").Value = "e2e\ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md"
$ws1.Range("E4").Value = "In Translation"
$ws1.Range("F4").Value = "In Translation"
$ws1.Range("G4").Value = "2016-08-29 22:42:50"

# Row 5 -> 1b48e175
$ws1.Range("A5").Value = "1b48e175-18de-4ddd-a054-b14e3ea4f762.md"
$ws1.Range("B5").Value = "e2e\1b48e175-18de-4ddd-a054-b14e3ea4f762.md"
$ws1.Range("E5").Value = "Ready for handoff"
$ws1.Range("F5").Value = "Ready for handoff"
$ws1.Range("G5").Value = "2016-08-29 22:41:39"

# Hyperlinks keep pointing at the same targets (rId2..rId5 / same URLs)
# as before the edit; only the display text follows the new row order.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb2f7ab57b14fd0b09522f42a75011eb4585b129/e2e/00f248cf-f57a-40f3-85ee-473dca7f5125.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\00f248cf-f57a-40f3-85ee-473dca7f5125.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/770d4aa91139fe93ceca103f217d3794b8b3e3fe/e2e/1b48e175-18de-4ddd-a054-b14e3ea4f762.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\3f14bde9-0c56-4a70-af96-9a02074d8301.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef167eee794e28926fcf6621deee466e564d0f01/e2e/3f14bde9-0c56-4a70-af96-9a02074d8301.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef167eee794e28926fcf6621deee466e564d0f01/e2e/ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md", [System.Type]::Missing, [System.Type]::Missing, "e2e\1b48e175-18de-4ddd-a054-b14e3ea4f762.md") | Out-Null

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

# Row 3 -> 3f14bde9
$ws2.Range("A3").Value = "3f14bde9-0c56-4a70-af96-9a02074d8301.md"
$ws2.Range("C3").Value = "In Translation"
$ws2.Range("G3").Value = "3f14bde9-0c56-4a70-af96-9a02074d8301.0d2737df5b05f8f9622cdb5ffc18a7ea43b5464b.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-29 22:42:45"

# Row 4 -> ba0dd1f4
$ws2.Range("A4").Value = "ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md"
$ws2.Range("C4").Value = "In Translation"
$ws2.Range("G4").Value = "ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.423cd28fb57506b2649e5eb8b503152ad7233acc.zh-cn.xlf"
$ws2.Range("H4").Value = "2016-08-29 22:42:45"

# Row 5 -> 1b48e175
$ws2.Range("A5").Value = "1b48e175-18de-4ddd-a054-b14e3ea4f762.md"
$ws2.Range("C5").Value = "Ready for handoff"
$ws2.Range("G5").Value = "1b48e175-18de-4ddd-a054-b14e3ea4f762.ad2dd3be08a37fc78801f7bb810ad29367702dec.zh-cn.xlf"
$ws2.Range("H5").Value = "2016-08-29 22:41:35"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb2f7ab57b14fd0b09522f42a75011eb4585b129/e2e/00f248cf-f57a-40f3-85ee-473dca7f5125.md", [System.Type]::Missing, [System.Type]::Missing, "00f248cf-f57a-40f3-85ee-473dca7f5125.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/95c30276a6a0f8bc52f0df4d1d2efad6f1181c13/e2e/00f248cf-f57a-40f3-85ee-473dca7f5125.md", [System.Type]::Missing, [System.Type]::Missing, "00f248cf-f57a-40f3-85ee-473dca7f5125.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/770d4aa91139fe93ceca103f217d3794b8b3e3fe/e2e/1b48e175-18de-4ddd-a054-b14e3ea4f762.md", [System.Type]::Missing, [System.Type]::Missing, "3f14bde9-0c56-4a70-af96-9a02074d8301.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef167eee794e28926fcf6621deee466e564d0f01/e2e/3f14bde9-0c56-4a70-af96-9a02074d8301.md", [System.Type]::Missing, [System.Type]::Missing, "ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef167eee794e28926fcf6621deee466e564d0f01/e2e/ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md", [System.Type]::Missing, [System.Type]::Missing, "1b48e175-18de-4ddd-a054-b14e3ea4f762.md") | Out-Null

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

# Row 3 -> 3f14bde9
$ws3.Range("A3").Value = "3f14bde9-0c56-4a70-af96-9a02074d8301.md"
$ws3.Range("C3").Value = "In Translation"
$ws3.Range("G3").Value = "3f14bde9-0c56-4a70-af96-9a02074d8301.0d2737df5b05f8f9622cdb5ffc18a7ea43b5464b.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-29 22:42:50"

# Row 4 -> ba0dd1f4
$ws3.Range("A4").Value = "ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md"
$ws3.Range("C4").Value = "In Translation"
$ws3.Range("G4").Value = "ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.423cd28fb57506b2649e5eb8b503152ad7233acc.de-de.xlf"
$ws3.Range("H4").Value = "2016-08-29 22:42:50"

# Row 5 -> 1b48e175
$ws3.Range("A5").Value = "1b48e175-18de-4ddd-a054-b14e3ea4f762.md"
$ws3.Range("C5").Value = "Ready for handoff"
$ws3.Range("G5").Value = "1b48e175-18de-4ddd-a054-b14e3ea4f762.ad2dd3be08a37fc78801f7bb810ad29367702dec.de-de.xlf"
$ws3.Range("H5").Value = "2016-08-29 22:41:39"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb2f7ab57b14fd0b09522f42a75011eb4585b129/e2e/00f248cf-f57a-40f3-85ee-473dca7f5125.md", [System.Type]::Missing, [System.Type]::Missing, "00f248cf-f57a-40f3-85ee-473dca7f5125.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6a1c57ba55a2de6cc9dddedb521b616a50b7ad29/e2e/00f248cf-f57a-40f3-85ee-473dca7f5125.md", [System.Type]::Missing, [System.Type]::Missing, "00f248cf-f57a-40f3-85ee-473dca7f5125.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/770d4aa91139fe93ceca103f217d3794b8b3e3fe/e2e/1b48e175-18de-4ddd-a054-b14e3ea4f762.md", [System.Type]::Missing, [System.Type]::Missing, "3f14bde9-0c56-4a70-af96-9a02074d8301.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef167eee794e28926fcf6621deee466e564d0f01/e2e/3f14bde9-0c56-4a70-af96-9a02074d8301.md", [System.Type]::Missing, [System.Type]::Missing, "ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ef167eee794e28926fcf6621deee466e564d0f01/e2e/ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md", [System.Type]::Missing, [System.Type]::Missing, "1b48e175-18de-4ddd-a054-b14e3ea4f762.md") | Out-Null
